$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.722.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.56%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.046.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -4.59%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.09%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'580.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.47%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'130.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.61%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.04%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.044.51"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -4.57%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.503"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.01%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -2.56%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'5.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.46%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.440"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.19%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000234"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.17%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'33.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.19%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +0.94%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.542.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.74%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'61.702.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.61%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.038.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.70%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.60%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'448.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.03%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -2.97%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.671"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -4.88%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -3.84%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'80.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.04%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'12.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -4.24%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.03%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.19%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -5.02%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.82%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'7.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -4.62%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'6.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -5.96%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'25.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -5.47%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.0973"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -6.41%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'2.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.29%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -6.29%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -3.25%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'50.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.59%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0₃0701"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.09%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.0373"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.68%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -1.22%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -2.68%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'380.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -5.58%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -6.90%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.697.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -5.42%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -0.02%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'123.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.44%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -4.34%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'34.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -7.56%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -6.36%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -2.84%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'23.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -6.69%  "
$ws.Range("E51").Style = "Normal"
